$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E3").Value = "Completed"
$ws.Range("D4").Value = "H001"
$ws.Range("E5").Value = "Completed"
$ws.Range("C6").Value = "H002"
$ws.Range("E6").Value = "Completed"

$ws.Range("F7").Select()
